$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 154 (shifts existing rows 154..186 -> 155..187,
# and Excel auto-extends the sheet's used range / dimension to A1:R187).
$ws.Rows.Item(154).Insert()

# Populate the new row 154. The surrounding rows in this table (Terminal La Palmera
# de La Serena / Ajo / Chino / Primera) share identical values for columns
# A, B, C, E, F, G, H, I, O, Q, R - only D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), N (Unidad de comercializacion)
# and P (Precio $/Kg) vary per record.
$ws.Range("A154").Value2 = 8
$ws.Range("B154").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C154").Value2 = "Coquimbo"
$ws.Range("D154").Value2 = 44543
$ws.Range("E154").Value2 = 4
$ws.Range("F154").Value2 = 100112003
$ws.Range("G154").Value2 = "Ajo"
$ws.Range("H154").Value2 = "Chino"
$ws.Range("I154").Value2 = "Primera"
$ws.Range("J154").Value2 = 500
$ws.Range("K154").Value2 = 19000
$ws.Range("L154").Value2 = 20000
$ws.Range("M154").Value2 = 19500
$ws.Range("N154").Value2 = "$/caja 10 kilos"
$ws.Range("O154").Value2 = "China"
$ws.Range("P154").Value2 = 1950
$ws.Range("Q154").Value2 = 10
$ws.Range("R154").Value2 = "Hortaliza"
